# Weekly price-sheet update for "Fruta, Terminal La Palmera de La Serena - Mango".
# A new week's worth of data (3 quality-grade rows: Especial / Primera / Segunda)
# is inserted at the top of the existing date-ordered block (row 390), pushing the
# rest of the block (and the dimension) down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 390, shifting existing rows 390:437 down to 393:440.
$ws.Rows("390:392").Insert()

# Shared (constant) column values for this product/market grouping.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100108
$producto    = "Tropicales y subtropicales"
$categoriaId = 100108002
$categoria   = "Mango"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 4 kilos"
$origen      = "Perú"
$kgUnidad    = 4

$fecha  = 44491
$volumen = 512
$precioMin = 6500
$precioMax = 7000
$precioProm = 6750
$precioKg = 1688

$calidades = @("Especial", "Primera", "Segunda")
for ($i = 0; $i -lt 3; $i++) {
    $r = 390 + $i
    $ws.Range("A$r").Value = $mercadoId
    $ws.Range("B$r").Value = $mercado
    $ws.Range("C$r").Value = $region
    $ws.Range("D$r").Value = $fecha
    $ws.Range("E$r").Value = $codreg
    $ws.Range("F$r").Value = $tipo
    $ws.Range("G$r").Value = $productoId
    $ws.Range("H$r").Value = $producto
    $ws.Range("I$r").Value = $categoriaId
    $ws.Range("J$r").Value = $categoria
    $ws.Range("K$r").Value = $variedad
    $ws.Range("L$r").Value = $calidades[$i]
    $ws.Range("M$r").Value = $volumen
    $ws.Range("N$r").Value = $precioMin
    $ws.Range("O$r").Value = $precioMax
    $ws.Range("P$r").Value = $precioProm
    $ws.Range("Q$r").Value = $unidad
    $ws.Range("R$r").Value = $origen
    $ws.Range("S$r").Value = $precioKg
    $ws.Range("T$r").Value = $kgUnidad
}
